$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 40
$ws.Range("A40").Value = 44621
$ws.Range("B40").Value = 0.4861111111111111
$ws.Range("C40").Value = 0.49583333333333335
$ws.Range("E40").Value = "Mail Chef de Projet"
$ws.Range("F40").Value = "mal entendu de ma part concernant la partie adaptative du pathfinding`nprévu initialement : un labyrinth: cela ne casse rien au niveau de ce qui a été mis en place, sa mise en place va impacter des modules qui sont encore à implémenter"

# Row 41
$ws.Range("A41").Value = 44621
$ws.Range("B41").Value = 0.49583333333333335
$ws.Range("C41").Value = 0.51041666666666663
$ws.Range("E41").Value = "Analyse et Conception`nPathfinding Labyrinth"
$ws.Range("F41").Value = "Problématique et intégration aux grid générics déjà en place"

# Row 42
$ws.Range("A42").Value = 44621
$ws.Range("B42").Value = 0.5625
$ws.Range("C42").Value = 0.57847222222222217
$ws.Range("E42").Value = "Redirection des objectifs"
$ws.Range("F42").Value = "La réponse au mail concernant la demande pour éviter la détection d'obstacle a été reçu.`nCependant ne restant que 3 semaines il est urgent de définir les objectifs afin d'avoir un projet qui remplisse un maximum le cahier des charges"
$ws.Range("G42").Value = "Mettre la documentation technique à jour"

# Row 43
$ws.Range("A43").Value = 44621
$ws.Range("B43").Value = 0.57847222222222217
$ws.Range("C43").Value = 0.62916666666666665
$ws.Range("E43").Value = "Analyse et conception`nPathfinding Detection chemin bloqué"
$ws.Range("H43").Value = "https://www.redblobgames.com/pathfinding/a-star/introduction.html"
$ws.Range("G43").Value = "Recherche via RedBlobGame contenant enormement d'article et lien d'articles sur le pathfinding"

# Row 44
$ws.Range("A44").Value = 44621
$ws.Range("B44").Value = 0.64236111111111105

# Row 40 grew to a two-line entry, so its height is no longer the old fixed 75pt
$ws.Rows.Item(40).RowHeight = 85.5

# Update selection to reflect last edited cell
$ws.Range("F44").Select()
